$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 606 (the "この子カンガルーは抱っこが好き" entry) entirely,
# shifting all subsequent rows up by one.
$ws.Rows("606:606").Delete()
